$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Total Compensation" (B) and "actual_partial_oh" (G) columns used to
# subtract four components (C-D-E-F / H-I-J-K). The extra subtractions meant
# duplicate asset-description rows (which doubled up E/F and J/K amounts)
# silently corrupted the partial-overhead total, so the formulas are
# simplified to only net the first two components.
#
# Rows 3:11 are written first as one fill so they land in a single shared
# formula group (matching how dragging the fill handle down a selection
# behaves), then row 2 is updated on its own afterwards.
$ws.Range("B3:B11").Formula = "=C3-D3"
$ws.Range("G3:G11").Formula = "=H3-I3"

$ws.Range("B2").Formula = "=C2-D2"
$ws.Range("G2").Formula = "=H2-I2"

# Move the active selection to G20, as recorded in the saved view state.
$ws.Range("G20").Select()
